# Refitting NCDEs to individual patients (for manuscript figure)
#
# Adds a "Label" column (H) flagging each row Control (0) / MDD (1),
# and refreshes the D/E/F refit values for the first (batchsize=100)
# block (rows 2-11) to the re-fit numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header "Label" in H1, matching the style of the other headers ---
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# --- Updated refit values (batchsize100 block: rows 2-8 and 11) ---
$ws.Range("D2").Value = 0.458367397482094
$ws.Range("E2").Value = 0.458367397482094

$ws.Range("D3").Value = 0.4062244361816142
$ws.Range("E3").Value = 0.4062244361816142

$ws.Range("D4").Value = 0.740189570332956
$ws.Range("E4").Value = 0.740189570332956

$ws.Range("D5").Value = 0.1799793559354867
$ws.Range("E5").Value = 0.1799793559354867

$ws.Range("D6").Value = 0.3462159766610102
$ws.Range("E6").Value = 0.3462159766610102

$ws.Range("D7").Value = 0.3567549837187902
$ws.Range("E7").Value = 0.6432450162812098

$ws.Range("D8").Value = 0.8016946051076239
$ws.Range("E8").Value = 0.1983053948923761

$ws.Range("D11").Value = 0.5341002434047126
$ws.Range("E11").Value = 0.4658997565952874
$ws.Range("F11").Value = 0.6263128519058228

# --- New "Label" column values: 0 for Control rows, 1 for MDD rows ---
# Block 1 (batchsize100): rows 2-6 = Control, rows 7-11 = MDD
$ws.Range("H2:H6").Value = 0
$ws.Range("H7:H11").Value = 1

# Block 2 (batchsize200): rows 12-16 = Control, rows 17-21 = MDD
$ws.Range("H12:H16").Value = 0
$ws.Range("H17:H21").Value = 1
